$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve header/date style (s="1") by copying an existing styled cell aside ---
$ws.Range("ZZ1").Value = 1
$ws.Range("B1").Copy($ws.Range("ZZ1"))

# --- Clear the whole used range to drop stale cells (old column BA, rows 23:24, etc.) ---
$ws.Range("A1:BA24").Clear()

# --- Re-apply the date/header style to the ranges that need it ---
$ws.Range("ZZ1").Copy($ws.Range("B1:AZ1"))
$ws.Range("ZZ1").Copy($ws.Range("A2:A22"))
$ws.Range("ZZ1").Clear()

# --- Populate header row (B1:AZ1) date values ---
$ws.Range("B1").Value = 39583
$ws.Range("C1").Value = 39765
$ws.Range("D1").Value = 39948
$ws.Range("E1").Value = 40130
$ws.Range("F1").Value = 40310
$ws.Range("G1").Value = 40494
$ws.Range("H1").Value = 40676
$ws.Range("I1").Value = 40862
$ws.Range("J1").Value = 41044
$ws.Range("K1").Value = 41228
$ws.Range("L1").Value = 41409
$ws.Range("M1").Value = 41592
$ws.Range("N1").Value = 41774
$ws.Range("O1").Value = 41957
$ws.Range("P1").Value = 42137
$ws.Range("Q1").Value = 42321
$ws.Range("R1").Value = 42503
$ws.Range("S1").Value = 42689
$ws.Range("T1").Value = 42867
$ws.Range("U1").Value = 43053
$ws.Range("V1").Value = 43145
$ws.Range("W1").Value = 43235
$ws.Range("X1").Value = 43326
$ws.Range("Y1").Value = 43418
$ws.Range("Z1").Value = 43510
$ws.Range("AA1").Value = 43600
$ws.Range("AB1").Value = 43691
$ws.Range("AC1").Value = 43783
$ws.Range("AD1").Value = 43875
$ws.Range("AE1").Value = 43966
$ws.Range("AF1").Value = 44068
$ws.Range("AG1").Value = 44159
$ws.Range("AH1").Value = 44251
$ws.Range("AI1").Value = 44341
$ws.Range("AJ1").Value = 44432
$ws.Range("AK1").Value = 44525
$ws.Range("AL1").Value = 44617
$ws.Range("AM1").Value = 44706
$ws.Range("AN1").Value = 44798
$ws.Range("AO1").Value = 44890
$ws.Range("AP1").Value = 44981
$ws.Range("AQ1").Value = 45071
$ws.Range("AR1").Value = 45163
$ws.Range("AS1").Value = 45254
$ws.Range("AT1").Value = 45345
$ws.Range("AU1").Value = 45436
$ws.Range("AV1").Value = 45534
$ws.Range("AW1").Value = 45618
$ws.Range("AX1").Value = 45713
$ws.Range("AY1").Value = 45800
$ws.Range("AZ1").Value = 45891
$ws.Range("A2").Value = 39813
$ws.Range("A3").Value = 40178
$ws.Range("E3").Value = -1.324983933426893
$ws.Range("F3").Value = -1.324983933426893
$ws.Range("G3").Value = -1.324983933426893
$ws.Range("H3").Value = -1.324983933426893
$ws.Range("I3").Value = -1.324983933426893
$ws.Range("J3").Value = -1.324983933426893
$ws.Range("K3").Value = -1.324983933426882
$ws.Range("L3").Value = -1.324983933426882
$ws.Range("M3").Value = -1.324983933426882
$ws.Range("N3").Value = -1.324983933426882
$ws.Range("O3").Value = -1.324983933426882
$ws.Range("P3").Value = -1.324983933426882
$ws.Range("Q3").Value = -1.324983933426882
$ws.Range("R3").Value = -1.324983933426882
$ws.Range("S3").Value = -1.324983933426882
$ws.Range("T3").Value = -1.324983933426882
$ws.Range("U3").Value = -1.324983933426882
$ws.Range("V3").Value = -1.324983933426882
$ws.Range("W3").Value = -1.324983933426882
$ws.Range("X3").Value = -1.324983933426882
$ws.Range("Y3").Value = -1.324983933426882
$ws.Range("Z3").Value = -1.324983933426882
$ws.Range("AA3").Value = -1.324983933426882
$ws.Range("AB3").Value = -1.324983933426882
$ws.Range("AC3").Value = -1.324983933426882
$ws.Range("AD3").Value = -1.324983933426882
$ws.Range("AE3").Value = -1.324983933426882
$ws.Range("AF3").Value = -1.324983933426882
$ws.Range("AG3").Value = -1.324983933426882
$ws.Range("AH3").Value = -1.324983933426882
$ws.Range("AI3").Value = -1.324983933426882
$ws.Range("AJ3").Value = -1.324983933426882
$ws.Range("AK3").Value = -1.324983933426882
$ws.Range("AL3").Value = -1.324983933426882
$ws.Range("AM3").Value = -1.324983933426882
$ws.Range("AN3").Value = -1.324983933426882
$ws.Range("AO3").Value = -1.324983933426882
$ws.Range("AP3").Value = -1.324983933426882
$ws.Range("AQ3").Value = -1.324983933426882
$ws.Range("AR3").Value = -1.324983933426882
$ws.Range("AS3").Value = -1.324983933426882
$ws.Range("AT3").Value = -1.324983933426882
$ws.Range("AU3").Value = -1.324983933426882
$ws.Range("AV3").Value = -1.324983933426882
$ws.Range("AW3").Value = -1.324983933426882
$ws.Range("AX3").Value = -1.324983933426882
$ws.Range("AY3").Value = -1.324983933426882
$ws.Range("AZ3").Value = -1.324983933426882
$ws.Range("A4").Value = 40543
$ws.Range("G4").Value = -0.3900454704678369
$ws.Range("H4").Value = -0.3900454704678369
$ws.Range("I4").Value = -0.3900454704678369
$ws.Range("J4").Value = -0.3900454704678369
$ws.Range("K4").Value = -0.3900454704678369
$ws.Range("L4").Value = -0.3900454704678369
$ws.Range("M4").Value = -0.3900454704678369
$ws.Range("N4").Value = -0.3900454704678369
$ws.Range("O4").Value = -0.3900454704678369
$ws.Range("P4").Value = -0.3900454704678369
$ws.Range("Q4").Value = -0.3900454704678369
$ws.Range("R4").Value = -0.3900454704678369
$ws.Range("S4").Value = -0.3900454704678369
$ws.Range("T4").Value = -0.3900454704678369
$ws.Range("U4").Value = -0.3900454704678369
$ws.Range("V4").Value = -0.3900454704678369
$ws.Range("W4").Value = -0.3900454704678369
$ws.Range("X4").Value = -0.3900454704678369
$ws.Range("Y4").Value = -0.3900454704678369
$ws.Range("Z4").Value = -0.3900454704678369
$ws.Range("AA4").Value = -0.3900454704678369
$ws.Range("AB4").Value = -0.3900454704678369
$ws.Range("AC4").Value = -0.3900454704678369
$ws.Range("AD4").Value = -0.3900454704678369
$ws.Range("AE4").Value = -0.3900454704678369
$ws.Range("AF4").Value = -0.3900454704678369
$ws.Range("AG4").Value = -0.3900454704678369
$ws.Range("AH4").Value = -0.3900454704678369
$ws.Range("AI4").Value = -0.3900454704678369
$ws.Range("AJ4").Value = -0.3900454704678369
$ws.Range("AK4").Value = -0.3900454704678369
$ws.Range("AL4").Value = -0.3900454704678369
$ws.Range("AM4").Value = -0.3900454704678369
$ws.Range("AN4").Value = -0.3900454704678369
$ws.Range("AO4").Value = -0.3900454704678369
$ws.Range("AP4").Value = -0.3900454704678369
$ws.Range("AQ4").Value = -0.3900454704678369
$ws.Range("AR4").Value = -0.3900454704678369
$ws.Range("AS4").Value = -0.3900454704678369
$ws.Range("AT4").Value = -0.3900454704678369
$ws.Range("AU4").Value = -0.3900454704678369
$ws.Range("AV4").Value = -0.3900454704678369
$ws.Range("AW4").Value = -0.3900454704678369
$ws.Range("AX4").Value = -0.3900454704678369
$ws.Range("AY4").Value = -0.3900454704678369
$ws.Range("AZ4").Value = -0.3900454704678369
$ws.Range("A5").Value = 40908
$ws.Range("I5").Value = -0.29958481534893
$ws.Range("J5").Value = -0.29958481534893
$ws.Range("K5").Value = -0.2995848153489522
$ws.Range("L5").Value = -0.2995848153489522
$ws.Range("M5").Value = -0.2995848153489522
$ws.Range("N5").Value = -0.2995848153489522
$ws.Range("O5").Value = -0.2995848153489522
$ws.Range("P5").Value = -0.2995848153489522
$ws.Range("Q5").Value = -0.2995848153489522
$ws.Range("R5").Value = -0.2995848153489522
$ws.Range("S5").Value = -0.2995848153489522
$ws.Range("T5").Value = -0.2995848153489522
$ws.Range("U5").Value = -0.2995848153489522
$ws.Range("V5").Value = -0.2995848153489522
$ws.Range("W5").Value = -0.2995848153489522
$ws.Range("X5").Value = -0.2995848153489522
$ws.Range("Y5").Value = -0.2995848153489522
$ws.Range("Z5").Value = -0.2995848153489522
$ws.Range("AA5").Value = -0.2995848153489522
$ws.Range("AB5").Value = -0.2995848153489522
$ws.Range("AC5").Value = -0.2995848153489522
$ws.Range("AD5").Value = -0.2995848153489522
$ws.Range("AE5").Value = -0.2995848153489522
$ws.Range("AF5").Value = -0.2995848153489522
$ws.Range("AG5").Value = -0.2995848153489522
$ws.Range("AH5").Value = -0.2995848153489522
$ws.Range("AI5").Value = -0.2995848153489522
$ws.Range("AJ5").Value = -0.2995848153489522
$ws.Range("AK5").Value = -0.2995848153489522
$ws.Range("AL5").Value = -0.2995848153489522
$ws.Range("AM5").Value = -0.2995848153489522
$ws.Range("AN5").Value = -0.2995848153489522
$ws.Range("AO5").Value = -0.2995848153489522
$ws.Range("AP5").Value = -0.2995848153489522
$ws.Range("AQ5").Value = -0.2995848153489522
$ws.Range("AR5").Value = -0.2995848153489522
$ws.Range("AS5").Value = -0.2995848153489522
$ws.Range("AT5").Value = -0.2995848153489522
$ws.Range("AU5").Value = -0.2995848153489522
$ws.Range("AV5").Value = -0.2995848153489522
$ws.Range("AW5").Value = -0.2995848153489522
$ws.Range("AX5").Value = -0.2995848153489522
$ws.Range("AY5").Value = -0.2995848153489522
$ws.Range("AZ5").Value = -0.2995848153489522
$ws.Range("A6").Value = 41274
$ws.Range("K6").Value = -0.2075757021743008
$ws.Range("L6").Value = -0.2075757021743008
$ws.Range("M6").Value = -0.2075757021743008
$ws.Range("N6").Value = -0.2075757021743008
$ws.Range("O6").Value = -0.2075757021743008
$ws.Range("P6").Value = -0.2075757021743008
$ws.Range("Q6").Value = -0.2075757021743008
$ws.Range("R6").Value = -0.2075757021743008
$ws.Range("S6").Value = -0.2075757021743008
$ws.Range("T6").Value = -0.2075757021743008
$ws.Range("U6").Value = -0.2075757021743008
$ws.Range("V6").Value = -0.2075757021743008
$ws.Range("W6").Value = -0.2075757021743008
$ws.Range("X6").Value = -0.2075757021743008
$ws.Range("Y6").Value = -0.2075757021743008
$ws.Range("Z6").Value = -0.2075757021743008
$ws.Range("AA6").Value = -0.2075757021743008
$ws.Range("AB6").Value = -0.2075757021743008
$ws.Range("AC6").Value = -0.2075757021743008
$ws.Range("AD6").Value = -0.2075757021743008
$ws.Range("AE6").Value = -0.2075757021743008
$ws.Range("AF6").Value = -0.2075757021743008
$ws.Range("AG6").Value = -0.2075757021743008
$ws.Range("AH6").Value = -0.2075757021743008
$ws.Range("AI6").Value = -0.2075757021743008
$ws.Range("AJ6").Value = -0.2075757021743008
$ws.Range("AK6").Value = -0.2075757021743008
$ws.Range("AL6").Value = -0.2075757021743008
$ws.Range("AM6").Value = -0.2075757021743008
$ws.Range("AN6").Value = -0.2075757021743008
$ws.Range("AO6").Value = -0.2075757021743008
$ws.Range("AP6").Value = -0.2075757021743008
$ws.Range("AQ6").Value = -0.2075757021743008
$ws.Range("AR6").Value = -0.2075757021743008
$ws.Range("AS6").Value = -0.2075757021743008
$ws.Range("AT6").Value = -0.2075757021743008
$ws.Range("AU6").Value = -0.2075757021743008
$ws.Range("AV6").Value = -0.2075757021743008
$ws.Range("AW6").Value = -0.2075757021743008
$ws.Range("AX6").Value = -0.2075757021743008
$ws.Range("AY6").Value = -0.2075757021743008
$ws.Range("AZ6").Value = -0.2075757021743008
$ws.Range("A7").Value = 41639
$ws.Range("K7").Value = -0.4223781730902543
$ws.Range("L7").Value = -0.06726240733891942
$ws.Range("M7").Value = 0.124712275190686
$ws.Range("N7").Value = 0.124712275190686
$ws.Range("O7").Value = 0.124712275190686
$ws.Range("P7").Value = 0.124712275190686
$ws.Range("Q7").Value = 0.124712275190686
$ws.Range("R7").Value = 0.124712275190686
$ws.Range("S7").Value = 0.124712275190686
$ws.Range("T7").Value = 0.124712275190686
$ws.Range("U7").Value = 0.124712275190686
$ws.Range("V7").Value = 0.124712275190686
$ws.Range("W7").Value = 0.124712275190686
$ws.Range("X7").Value = 0.124712275190686
$ws.Range("Y7").Value = 0.124712275190686
$ws.Range("Z7").Value = 0.124712275190686
$ws.Range("AA7").Value = 0.124712275190686
$ws.Range("AB7").Value = 0.124712275190686
$ws.Range("AC7").Value = 0.124712275190686
$ws.Range("AD7").Value = 0.124712275190686
$ws.Range("AE7").Value = 0.124712275190686
$ws.Range("AF7").Value = 0.124712275190686
$ws.Range("AG7").Value = 0.124712275190686
$ws.Range("AH7").Value = 0.124712275190686
$ws.Range("AI7").Value = 0.124712275190686
$ws.Range("AJ7").Value = 0.124712275190686
$ws.Range("AK7").Value = 0.124712275190686
$ws.Range("AL7").Value = 0.124712275190686
$ws.Range("AM7").Value = 0.124712275190686
$ws.Range("AN7").Value = 0.124712275190686
$ws.Range("AO7").Value = 0.124712275190686
$ws.Range("AP7").Value = 0.124712275190686
$ws.Range("AQ7").Value = 0.124712275190686
$ws.Range("AR7").Value = 0.124712275190686
$ws.Range("AS7").Value = 0.124712275190686
$ws.Range("AT7").Value = 0.124712275190686
$ws.Range("AU7").Value = 0.124712275190686
$ws.Range("AV7").Value = 0.124712275190686
$ws.Range("AW7").Value = 0.124712275190686
$ws.Range("AX7").Value = 0.124712275190686
$ws.Range("AY7").Value = 0.124712275190686
$ws.Range("AZ7").Value = 0.124712275190686
$ws.Range("A8").Value = 42004
$ws.Range("K8").Value = -0.5490886506258952
$ws.Range("L8").Value = -0.4774698422615242
$ws.Range("M8").Value = -0.3577371449824729
$ws.Range("N8").Value = -0.3584575688954428
$ws.Range("O8").Value = -0.255298189276465
$ws.Range("P8").Value = -0.255298189276465
$ws.Range("Q8").Value = -0.255298189276465
$ws.Range("R8").Value = -0.255298189276465
$ws.Range("S8").Value = -0.255298189276465
$ws.Range("T8").Value = -0.255298189276465
$ws.Range("U8").Value = -0.255298189276465
$ws.Range("V8").Value = -0.255298189276465
$ws.Range("W8").Value = -0.255298189276465
$ws.Range("X8").Value = -0.255298189276465
$ws.Range("Y8").Value = -0.255298189276465
$ws.Range("Z8").Value = -0.255298189276465
$ws.Range("AA8").Value = -0.255298189276465
$ws.Range("AB8").Value = -0.255298189276465
$ws.Range("AC8").Value = -0.255298189276465
$ws.Range("AD8").Value = -0.255298189276465
$ws.Range("AE8").Value = -0.255298189276465
$ws.Range("AF8").Value = -0.255298189276465
$ws.Range("AG8").Value = -0.255298189276465
$ws.Range("AH8").Value = -0.255298189276465
$ws.Range("AI8").Value = -0.255298189276465
$ws.Range("AJ8").Value = -0.255298189276465
$ws.Range("AK8").Value = -0.255298189276465
$ws.Range("AL8").Value = -0.255298189276465
$ws.Range("AM8").Value = -0.255298189276465
$ws.Range("AN8").Value = -0.255298189276465
$ws.Range("AO8").Value = -0.255298189276465
$ws.Range("AP8").Value = -0.255298189276465
$ws.Range("AQ8").Value = -0.255298189276465
$ws.Range("AR8").Value = -0.255298189276465
$ws.Range("AS8").Value = -0.255298189276465
$ws.Range("AT8").Value = -0.255298189276465
$ws.Range("AU8").Value = -0.255298189276465
$ws.Range("AV8").Value = -0.255298189276465
$ws.Range("AW8").Value = -0.255298189276465
$ws.Range("AX8").Value = -0.255298189276465
$ws.Range("AY8").Value = -0.255298189276465
$ws.Range("AZ8").Value = -0.255298189276465
$ws.Range("A9").Value = 42369
$ws.Range("L9").Value = -0.4828982001363724
$ws.Range("M9").Value = -0.4292407320315994
$ws.Range("N9").Value = -0.4112436562971
$ws.Range("O9").Value = -0.2877346565283379
$ws.Range("P9").Value = -0.00968885111266582
$ws.Range("Q9").Value = 0.07418514192796266
$ws.Range("R9").Value = 0.07418514192796266
$ws.Range("S9").Value = 0.07418514192796266
$ws.Range("T9").Value = 0.07418514192796266
$ws.Range("U9").Value = 0.07418514192796266
$ws.Range("V9").Value = 0.07418514192796266
$ws.Range("W9").Value = 0.07418514192796266
$ws.Range("X9").Value = 0.07418514192796266
$ws.Range("Y9").Value = 0.07418514192796266
$ws.Range("Z9").Value = 0.07418514192796266
$ws.Range("AA9").Value = 0.07418514192796266
$ws.Range("AB9").Value = 0.07418514192796266
$ws.Range("AC9").Value = 0.07418514192796266
$ws.Range("AD9").Value = 0.07418514192796266
$ws.Range("AE9").Value = 0.07418514192796266
$ws.Range("AF9").Value = 0.07418514192796266
$ws.Range("AG9").Value = 0.07418514192796266
$ws.Range("AH9").Value = 0.07418514192796266
$ws.Range("AI9").Value = 0.07418514192796266
$ws.Range("AJ9").Value = 0.07418514192796266
$ws.Range("AK9").Value = 0.07418514192796266
$ws.Range("AL9").Value = 0.07418514192796266
$ws.Range("AM9").Value = 0.07418514192796266
$ws.Range("AN9").Value = 0.07418514192796266
$ws.Range("AO9").Value = 0.07418514192796266
$ws.Range("AP9").Value = 0.07418514192796266
$ws.Range("AQ9").Value = 0.07418514192796266
$ws.Range("AR9").Value = 0.07418514192796266
$ws.Range("AS9").Value = 0.07418514192796266
$ws.Range("AT9").Value = 0.07418514192796266
$ws.Range("AU9").Value = 0.07418514192796266
$ws.Range("AV9").Value = 0.07418514192796266
$ws.Range("AW9").Value = 0.07418514192796266
$ws.Range("AX9").Value = 0.07418514192796266
$ws.Range("AY9").Value = 0.07418514192796266
$ws.Range("AZ9").Value = 0.07418514192796266
$ws.Range("A10").Value = 42735
$ws.Range("N10").Value = -0.4304242973383055
$ws.Range("O10").Value = -0.3979953014972226
$ws.Range("P10").Value = -0.3537462851234685
$ws.Range("Q10").Value = -0.2097319935285391
$ws.Range("R10").Value = -0.1132037832954791
$ws.Range("S10").Value = -0.07611406013281474
$ws.Range("T10").Value = -0.07611406013281474
$ws.Range("U10").Value = -0.07611406013281474
$ws.Range("V10").Value = -0.07611406013281474
$ws.Range("W10").Value = -0.07611406013281474
$ws.Range("X10").Value = -0.07611406013281474
$ws.Range("Y10").Value = -0.07611406013281474
$ws.Range("Z10").Value = -0.07611406013281474
$ws.Range("AA10").Value = -0.07611406013281474
$ws.Range("AB10").Value = -0.07611406013281474
$ws.Range("AC10").Value = -0.07611406013281474
$ws.Range("AD10").Value = -0.07611406013281474
$ws.Range("AE10").Value = -0.07611406013281474
$ws.Range("AF10").Value = -0.07611406013281474
$ws.Range("AG10").Value = -0.07611406013281474
$ws.Range("AH10").Value = -0.07611406013281474
$ws.Range("AI10").Value = -0.07611406013281474
$ws.Range("AJ10").Value = -0.07611406013281474
$ws.Range("AK10").Value = -0.07611406013281474
$ws.Range("AL10").Value = -0.07611406013281474
$ws.Range("AM10").Value = -0.07611406013281474
$ws.Range("AN10").Value = -0.07611406013281474
$ws.Range("AO10").Value = -0.07611406013281474
$ws.Range("AP10").Value = -0.07611406013281474
$ws.Range("AQ10").Value = -0.07611406013281474
$ws.Range("AR10").Value = -0.07611406013281474
$ws.Range("AS10").Value = -0.07611406013281474
$ws.Range("AT10").Value = -0.07611406013281474
$ws.Range("AU10").Value = -0.07611406013281474
$ws.Range("AV10").Value = -0.07611406013281474
$ws.Range("AW10").Value = -0.07611406013281474
$ws.Range("AX10").Value = -0.07611406013281474
$ws.Range("AY10").Value = -0.07611406013281474
$ws.Range("AZ10").Value = -0.07611406013281474
$ws.Range("A11").Value = 43100
$ws.Range("P11").Value = -0.3602216788231694
$ws.Range("Q11").Value = -0.3346105473710614
$ws.Range("R11").Value = -0.3204027102583273
$ws.Range("S11").Value = -0.3036570471216304
$ws.Range("T11").Value = -0.2638577853126156
$ws.Range("U11").Value = -0.191300579729714
$ws.Range("V11").Value = -0.191300579729714
$ws.Range("W11").Value = -0.191300579729714
$ws.Range("X11").Value = -0.191300579729714
$ws.Range("Y11").Value = -0.191300579729714
$ws.Range("Z11").Value = -0.191300579729714
$ws.Range("AA11").Value = -0.191300579729714
$ws.Range("AB11").Value = -0.191300579729714
$ws.Range("AC11").Value = -0.191300579729714
$ws.Range("AD11").Value = -0.191300579729714
$ws.Range("AE11").Value = -0.191300579729714
$ws.Range("AF11").Value = -0.191300579729714
$ws.Range("AG11").Value = -0.191300579729714
$ws.Range("AH11").Value = -0.191300579729714
$ws.Range("AI11").Value = -0.191300579729714
$ws.Range("AJ11").Value = -0.191300579729714
$ws.Range("AK11").Value = -0.191300579729714
$ws.Range("AL11").Value = -0.191300579729714
$ws.Range("AM11").Value = -0.191300579729714
$ws.Range("AN11").Value = -0.191300579729714
$ws.Range("AO11").Value = -0.191300579729714
$ws.Range("AP11").Value = -0.191300579729714
$ws.Range("AQ11").Value = -0.191300579729714
$ws.Range("AR11").Value = -0.191300579729714
$ws.Range("AS11").Value = -0.191300579729714
$ws.Range("AT11").Value = -0.191300579729714
$ws.Range("AU11").Value = -0.191300579729714
$ws.Range("AV11").Value = -0.191300579729714
$ws.Range("AW11").Value = -0.191300579729714
$ws.Range("AX11").Value = -0.191300579729714
$ws.Range("AY11").Value = -0.191300579729714
$ws.Range("AZ11").Value = -0.191300579729714
$ws.Range("A12").Value = 43465
$ws.Range("R12").Value = -0.3211834654844647
$ws.Range("S12").Value = -0.3142097218752804
$ws.Range("T12").Value = -0.3143428957755656
$ws.Range("U12").Value = -0.1907914728172644
$ws.Range("V12").Value = -0.05059932794906352
$ws.Range("W12").Value = 0.04166709579394023
$ws.Range("X12").Value = 0.1246424644191668
$ws.Range("Y12").Value = 0.0970330232288763
$ws.Range("Z12").Value = 0.0970330232288763
$ws.Range("AA12").Value = 0.0970330232288763
$ws.Range("AB12").Value = 0.0970330232288763
$ws.Range("AC12").Value = 0.0970330232288763
$ws.Range("AD12").Value = 0.0970330232288763
$ws.Range("AE12").Value = 0.0970330232288763
$ws.Range("AF12").Value = 0.0970330232288763
$ws.Range("AG12").Value = 0.0970330232288763
$ws.Range("AH12").Value = 0.0970330232288763
$ws.Range("AI12").Value = 0.0970330232288763
$ws.Range("AJ12").Value = 0.0970330232288763
$ws.Range("AK12").Value = 0.0970330232288763
$ws.Range("AL12").Value = 0.0970330232288763
$ws.Range("AM12").Value = 0.0970330232288763
$ws.Range("AN12").Value = 0.0970330232288763
$ws.Range("AO12").Value = 0.0970330232288763
$ws.Range("AP12").Value = 0.0970330232288763
$ws.Range("AQ12").Value = 0.0970330232288763
$ws.Range("AR12").Value = 0.0970330232288763
$ws.Range("AS12").Value = 0.0970330232288763
$ws.Range("AT12").Value = 0.0970330232288763
$ws.Range("AU12").Value = 0.0970330232288763
$ws.Range("AV12").Value = 0.0970330232288763
$ws.Range("AW12").Value = 0.0970330232288763
$ws.Range("AX12").Value = 0.0970330232288763
$ws.Range("AY12").Value = 0.0970330232288763
$ws.Range("AZ12").Value = 0.0970330232288763
$ws.Range("A13").Value = 43830
$ws.Range("T13").Value = -0.3114748837725667
$ws.Range("U13").Value = -0.293687056619063
$ws.Range("V13").Value = -0.2866918525404771
$ws.Range("W13").Value = -0.2532926704812977
$ws.Range("X13").Value = -0.1729167731976178
$ws.Range("Y13").Value = -0.2529765062333933
$ws.Range("Z13").Value = -0.567279386573194
$ws.Range("AA13").Value = -0.7671134292608239
$ws.Range("AB13").Value = -0.6865797156070164
$ws.Range("AC13").Value = -0.7407518902333265
$ws.Range("AD13").Value = -0.7407518902333265
$ws.Range("AE13").Value = -0.7407518902333265
$ws.Range("AF13").Value = -0.7407518902333265
$ws.Range("AG13").Value = -0.7407518902333265
$ws.Range("AH13").Value = -0.7407518902333265
$ws.Range("AI13").Value = -0.7407518902333265
$ws.Range("AJ13").Value = -0.7407518902333265
$ws.Range("AK13").Value = -0.7407518902333265
$ws.Range("AL13").Value = -0.7407518902333265
$ws.Range("AM13").Value = -0.7407518902333265
$ws.Range("AN13").Value = -0.7407518902333265
$ws.Range("AO13").Value = -0.7407518902333265
$ws.Range("AP13").Value = -0.7407518902333265
$ws.Range("AQ13").Value = -0.7407518902333265
$ws.Range("AR13").Value = -0.7407518902333265
$ws.Range("AS13").Value = -0.7407518902333265
$ws.Range("AT13").Value = -0.7407518902333265
$ws.Range("AU13").Value = -0.7407518902333265
$ws.Range("AV13").Value = -0.7407518902333265
$ws.Range("AW13").Value = -0.7407518902333265
$ws.Range("AX13").Value = -0.7407518902333265
$ws.Range("AY13").Value = -0.7407518902333265
$ws.Range("AZ13").Value = -0.7407518902333265
$ws.Range("A14").Value = 44196
$ws.Range("W14").Value = -0.2681128372844399
$ws.Range("X14").Value = -0.2532122198970588
$ws.Range("Y14").Value = -0.2609030463090245
$ws.Range("Z14").Value = -0.2733008487008526
$ws.Range("AA14").Value = -0.3480983700859808
$ws.Range("AB14").Value = -0.2709243735070865
$ws.Range("AC14").Value = -0.428077259747528
$ws.Range("AD14").Value = -0.6415172061831176
$ws.Range("AE14").Value = -0.2043373675692961
$ws.Range("AF14").Value = 0.3056679541520335
$ws.Range("AG14").Value = 0.3056679541520335
$ws.Range("AH14").Value = 0.3056679541520335
$ws.Range("AI14").Value = 0.3056679541520335
$ws.Range("AJ14").Value = 0.3056679541520335
$ws.Range("AK14").Value = 0.3056679541520335
$ws.Range("AL14").Value = 0.3056679541520335
$ws.Range("AM14").Value = 0.3056679541520335
$ws.Range("AN14").Value = 0.3056679541520335
$ws.Range("AO14").Value = 0.3056679541520335
$ws.Range("AP14").Value = 0.3056679541520335
$ws.Range("AQ14").Value = 0.3056679541520335
$ws.Range("AR14").Value = 0.3056679541520335
$ws.Range("AS14").Value = 0.3056679541520335
$ws.Range("AT14").Value = 0.3056679541520335
$ws.Range("AU14").Value = 0.3056679541520335
$ws.Range("AV14").Value = 0.3056679541520335
$ws.Range("AW14").Value = 0.3056679541520335
$ws.Range("AX14").Value = 0.3056679541520335
$ws.Range("AY14").Value = 0.3056679541520335
$ws.Range("AZ14").Value = 0.3056679541520335
$ws.Range("A15").Value = 44561
$ws.Range("AA15").Value = -0.3127278713534576
$ws.Range("AB15").Value = -0.2985160331797965
$ws.Range("AC15").Value = -0.3128324057260823
$ws.Range("AD15").Value = -0.3202690301181033
$ws.Range("AE15").Value = -0.1872494743064723
$ws.Range("AF15").Value = -0.02738999065564629
$ws.Range("AG15").Value = -0.2112001730687485
$ws.Range("AH15").Value = -0.6729198470149966
$ws.Range("AI15").Value = -0.8080927309597863
$ws.Range("AJ15").Value = -1.388491535160907
$ws.Range("AK15").Value = -1.388491535160907
$ws.Range("AL15").Value = -1.388491535160907
$ws.Range("AM15").Value = -1.388491535160907
$ws.Range("AN15").Value = -1.388491535160907
$ws.Range("AO15").Value = -1.388491535160907
$ws.Range("AP15").Value = -1.388491535160907
$ws.Range("AQ15").Value = -1.388491535160907
$ws.Range("AR15").Value = -1.388491535160907
$ws.Range("AS15").Value = -1.388491535160907
$ws.Range("AT15").Value = -1.388491535160907
$ws.Range("AU15").Value = -1.388491535160907
$ws.Range("AV15").Value = -1.388491535160907
$ws.Range("AW15").Value = -1.388491535160907
$ws.Range("AX15").Value = -1.388491535160907
$ws.Range("AY15").Value = -1.388491535160907
$ws.Range("AZ15").Value = -1.388491535160907
$ws.Range("A16").Value = 44926
$ws.Range("AE16").Value = -0.2863691763874465
$ws.Range("AF16").Value = -0.1764633559520723
$ws.Range("AG16").Value = -0.1174713876879729
$ws.Range("AH16").Value = -0.1659341199724107
$ws.Range("AI16").Value = -0.1644030883838465
$ws.Range("AJ16").Value = -2.365218846580541
$ws.Range("AK16").Value = -1.701252732314051
$ws.Range("AL16").Value = -1.748780309677478
$ws.Range("AM16").Value = -1.697148566375528
$ws.Range("AN16").Value = -1.678482969789596
$ws.Range("AO16").Value = -1.678482969789596
$ws.Range("AP16").Value = -1.678482969789596
$ws.Range("AQ16").Value = -1.678482969789596
$ws.Range("AR16").Value = -1.678482969789596
$ws.Range("AS16").Value = -1.678482969789596
$ws.Range("AT16").Value = -1.678482969789596
$ws.Range("AU16").Value = -1.678482969789596
$ws.Range("AV16").Value = -1.678482969789596
$ws.Range("AW16").Value = -1.678482969789596
$ws.Range("AX16").Value = -1.678482969789596
$ws.Range("AY16").Value = -1.678482969789596
$ws.Range("AZ16").Value = -1.678482969789596
$ws.Range("A17").Value = 45291
$ws.Range("AH17").Value = -0.2262453279458954
$ws.Range("AI17").Value = -0.1915283081898367
$ws.Range("AJ17").Value = -0.2395344214953798
$ws.Range("AK17").Value = -0.3039708008245712
$ws.Range("AL17").Value = -0.3138757159763084
$ws.Range("AM17").Value = -0.2339923140600275
$ws.Range("AN17").Value = -0.1583814731481836
$ws.Range("AO17").Value = -0.6241481568271312
$ws.Range("AP17").Value = -0.6410531931378527
$ws.Range("AQ17").Value = -0.6247846736575413
$ws.Range("AR17").Value = -0.5999457276250508
$ws.Range("AS17").Value = -0.5999457276250508
$ws.Range("AT17").Value = -0.5999457276250508
$ws.Range("AU17").Value = -0.5999457276250508
$ws.Range("AV17").Value = -0.5999457276250508
$ws.Range("AW17").Value = -0.5999457276250508
$ws.Range("AX17").Value = -0.5999457276250508
$ws.Range("AY17").Value = -0.5999457276250508
$ws.Range("AZ17").Value = -0.5999457276250508
$ws.Range("A18").Value = 45657
$ws.Range("AL18").Value = -0.3686474460577349
$ws.Range("AM18").Value = -0.3294995785474542
$ws.Range("AN18").Value = -0.3073589516970454
$ws.Range("AO18").Value = -0.3568684571922609
$ws.Range("AP18").Value = -0.3624456929795739
$ws.Range("AQ18").Value = -0.3554771869619944
$ws.Range("AR18").Value = -0.2733934973402352
$ws.Range("AS18").Value = -0.2119687890143274
$ws.Range("AT18").Value = -0.1344629014043419
$ws.Range("AU18").Value = -0.09040308684795662
$ws.Range("AV18").Value = -0.05499271238530445
$ws.Range("AW18").Value = -0.05499271238530445
$ws.Range("AX18").Value = -0.05499271238530445
$ws.Range("AY18").Value = -0.05499271238530445
$ws.Range("AZ18").Value = -0.05499271238530445
$ws.Range("A19").Value = 46022
$ws.Range("AP19").Value = -0.3726310378213471
$ws.Range("AQ19").Value = -0.3811678132403018
$ws.Range("AR19").Value = -0.3780242922104993
$ws.Range("AS19").Value = -0.3675939520929039
$ws.Range("AT19").Value = -0.356174448876545
$ws.Range("AU19").Value = -0.3159097170635006
$ws.Range("AV19").Value = -0.1967532854572851
$ws.Range("AW19").Value = -0.06564014165270082
$ws.Range("AX19").Value = 0.393650197209916
$ws.Range("AY19").Value = 0.3813481955213138
$ws.Range("AZ19").Value = 0.4335297397760618
$ws.Range("A20").Value = 46387
$ws.Range("AT20").Value = -0.370253819370725
$ws.Range("AU20").Value = -0.3674463749899881
$ws.Range("AV20").Value = -0.3611597286974577
$ws.Range("AW20").Value = -0.3724803060494719
$ws.Range("AX20").Value = -0.2112949500827632
$ws.Range("AY20").Value = -0.2523778956734835
$ws.Range("AZ20").Value = -0.05771130674934177
$ws.Range("A21").Value = 46752
$ws.Range("AX21").Value = -0.3179113751378249
$ws.Range("AY21").Value = -0.3521544584122904
$ws.Range("AZ21").Value = -0.319697601342106
$ws.Range("A22").Value = 47118